# ToDoList_Form.xlsx edit script
# - Team sheet: recategorize rows 2-7 under "바디 전실", rename a couple of
#   entries, add a "분류No" column, and drop the old scratch rows (8-24).
# - Person sheet: rename a couple of entries, add a "분류No" column, and
#   drop the old trailing scratch row (11).
# - DB sheet: remove the old "교육 예산 관리" task row (everything below
#   shifts up one row), renumber the affected task IDs in column A, and
#   append a new "자동화 프로그램 개발" task row at the bottom.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Team sheet
# ---------------------------------------------------------------------
$wsTeam = $wb.Worksheets.Item("Team")

$wsTeam.Range("A8:A24").EntireRow.Delete()

$wsTeam.Range("C1").Value = "분류No"

$wsTeam.Range("A2").Value = "바디 전실"
$wsTeam.Range("A3").Value = "바디 전실"
$wsTeam.Range("A4").Value = "바디 전실"
$wsTeam.Range("A5").Value = "바디 전실"
$wsTeam.Range("A6").Value = "바디 전실"
$wsTeam.Range("A7").Value = "바디 전실"

$wsTeam.Range("B6").Value = "내장실"
$wsTeam.Range("B7").Value = "안전실"

$wsTeam.Range("C2").Value = 1
$wsTeam.Range("C3").Value = 1
$wsTeam.Range("C4").Value = 1
$wsTeam.Range("C5").Value = 1
$wsTeam.Range("C6").Value = 1
$wsTeam.Range("C7").Value = 1

# ---------------------------------------------------------------------
# Person sheet
# ---------------------------------------------------------------------
$wsPerson = $wb.Worksheets.Item("Person")

$wsPerson.Range("A11:A11").EntireRow.Delete()

$wsPerson.Range("C1").Value = "분류No"

$wsPerson.Range("B3").Value = "외장책임"
$wsPerson.Range("B4").Value = "내장책임"

$wsPerson.Range("A7").Value = "원가"
$wsPerson.Range("B7").Value = "가가가"
$wsPerson.Range("B8").Value = "나나나"
$wsPerson.Range("B9").Value = "다다다"
$wsPerson.Range("B10").Value = "라라라"

$wsPerson.Range("C2").Value = 1
$wsPerson.Range("C3").Value = 1
$wsPerson.Range("C4").Value = 1
$wsPerson.Range("C5").Value = 2
$wsPerson.Range("C6").Value = 2
$wsPerson.Range("C7").Value = 3
$wsPerson.Range("C8").Value = 3
$wsPerson.Range("C9").Value = 3
$wsPerson.Range("C10").Value = 3

# ---------------------------------------------------------------------
# DB sheet
# ---------------------------------------------------------------------
$wsDB = $wb.Worksheets.Item("DB")

# Drop the old "교육 예산 관리" row; rows 19-50 shift up to 18-49.
$wsDB.Range("A18:A18").EntireRow.Delete()

# Renumber the task IDs in column A for the shifted rows.
$wsDB.Range("A18").Value = "03-00-00"
$wsDB.Range("A19").Value = "03-01-00"
$wsDB.Range("A20").Value = "03-01-01"
$wsDB.Range("A21").Value = "03-01-02"
$wsDB.Range("A22").Value = "03-01-03"
$wsDB.Range("A23").Value = "03-02-00"
$wsDB.Range("A24").Value = "03-02-01"
$wsDB.Range("A25").Value = "03-02-02"
$wsDB.Range("A26").Value = "03-02-03"
$wsDB.Range("A27").Value = "03-02-04"
$wsDB.Range("A28").Value = "03-02-05"
$wsDB.Range("A29").Value = "03-02-06"
$wsDB.Range("A30").Value = "03-02-07"
$wsDB.Range("A31").Value = "03-02-08"
$wsDB.Range("A32").Value = "03-02-09"
$wsDB.Range("A33").Value = "03-03-00"
$wsDB.Range("A34").Value = "04-00-00"
$wsDB.Range("A35").Value = "04-01-00"
$wsDB.Range("A36").Value = "04-02-00"
$wsDB.Range("A37").Value = "04-03-00"
$wsDB.Range("A38").Value = "04-04-00"
$wsDB.Range("A39").Value = "05-00-00"
$wsDB.Range("A40").Value = "06-00-00"
$wsDB.Range("A41").Value = "06-01-00"
$wsDB.Range("A42").Value = "06-02-00"
$wsDB.Range("A43").Value = "06-03-00"
$wsDB.Range("A44").Value = "06-04-00"
$wsDB.Range("A45").Value = "06-05-00"
$wsDB.Range("A46").Value = "06-06-00"
$wsDB.Range("A47").Value = "06-07-00"
$wsDB.Range("A48").Value = "06-08-00"
$wsDB.Range("A49").Value = "06-09-00"

# Append the brand-new task row at the bottom.
$wsDB.Range("A50").Value = "07-00-00"
$wsDB.Range("B50").Value = "자동화 프로그램 개발"
$wsDB.Range("C50").Value = " "
$wsDB.Range("D50").Value = " "
$wsDB.Range("E50").Value = "진행"
$wsDB.Range("F50").Value = " "
$wsDB.Range("G50").Value = " "
$wsDB.Range("H50").Value = " "
$wsDB.Range("I50").Value = 1

Write-Output "edit complete"
